# Soft margin SVM with and without PCA - fill in results rows 17 (Normal) and 18 (PCA)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Data")

# Row 17 - Soft SVM Normal
$ws.Range("B17").Value = 0.96613899999999997
$ws.Range("C17").Value = 0.96605300000000005
$ws.Range("D17").Value = 0.965144
$ws.Range("E17").Value = 0.96391700000000002
$ws.Range("F17").Value = 0.95434399999999997
$ws.Range("G17").Value = 0.95571899999999999

$ws.Range("K17").Value = 0.92058799999999996
$ws.Range("L17").Value = 0.91176400000000002
$ws.Range("M17").Value = 0.90882300000000005
$ws.Range("N17").Value = 0.89705800000000002
$ws.Range("O17").Value = 0.88235200000000003
$ws.Range("P17").Value = 0.86176399999999997

$ws.Range("R17").Value = "Tim"

# Row 18 - Soft SVM PCA
$ws.Range("B18").Value = 0.95362499999999994
$ws.Range("C18").Value = 0.95460100000000003
$ws.Range("D18").Value = 0.95238
$ws.Range("E18").Value = 0.94992600000000005
$ws.Range("F18").Value = 0.94403499999999996
$ws.Range("G18").Value = 0.96309900000000004

$ws.Range("K18").Value = 0.91470499999999999
$ws.Range("L18").Value = 0.90294099999999999
$ws.Range("M18").Value = 0.90588199999999997
$ws.Range("N18").Value = 0.87646999999999997
$ws.Range("O18").Value = 0.888235
$ws.Range("P18").Value = 0.83529399999999998

$ws.Range("R18").Value = "Tim"

# Update the selected cell to mirror the author's last cursor position
$ws.Range("O28").Select()

$wb.Save()
